$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '40.102.81'
$ws.Range('E2').Value = '  -2.40%  '
$ws.Range('D3').Value = '2.343.51'
$ws.Range('E3').Value = '  -3.51%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range("D5").Value = "'310.90"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.59%  '
$ws.Range("D6").Value = "'85.64"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.95%  '
$ws.Range("D7").Value = "'0.531"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.76%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range("D9").Value = "'0.485"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.92%  '
$ws.Range("D10").Value = "'0.0811"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.65%  '
$ws.Range("D11").Value = "'30.08"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -6.06%  '
$ws.Range("D12").Value = "'0.110"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.28%  '
$ws.Range('D13').Value = '2.702.08'
$ws.Range('E13').Value = '  -3.97%  '
$ws.Range('E14').Value = '  -3.94%  '
$ws.Range("D15").Value = "'14.76"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.90%  '
$ws.Range('D16').Value = '2.369.94'
$ws.Range('E16').Value = '  -2.88%  '
$ws.Range("D17").Value = "'0.760"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.10%  '
$ws.Range('D18').Value = '40.062.05'
$ws.Range('E18').Value = '  -2.49%  '
$ws.Range('D19').Value = '0.0₃0905'
$ws.Range('E19').Value = '  -1.86%  '
$ws.Range('E20').Value = '  -2.03%  '
$ws.Range("D21").Value = "'68.23"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.39%  '
$ws.Range("D22").Value = "'10.71"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.58%  '
$ws.Range("D23").Value = "'235.69"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('E24').Value = '  -4.90%  '
$ws.Range('E26').Value = '  -2.79%  '
$ws.Range("D27").Value = "'23.38"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.20%  '
$ws.Range("D28").Value = "'2.13"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.67%  '
$ws.Range("D29").Value = "'9.30"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.28%  '
$ws.Range("D30").Value = "'34.74"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range("D31").Value = "'153.29"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.46%  '
$ws.Range('E32').Value = '  -0.10%  '
$ws.Range("D33").Value = "'5.12"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.38%  '
$ws.Range("D34").Value = "'2.45"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.72%  '
$ws.Range("D35").Value = "'0.0721"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.03%  '
$ws.Range('E36').Value = '  -0.57%  '
$ws.Range("D37").Value = "'2.83"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.35%  '
$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D38").Value = "'1.73"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.18%  '
$ws.Range('B39').Value = 'Celestia'
$ws.Range('C39').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D39").Value = "'15.67"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -5.73%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").Value = "'0.0986"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.53%  '
$ws.Range("D41").Value = "'3.88"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.80%  '
$ws.Range('D42').Value = '1.958.71'
$ws.Range('E42').Value = '  -1.39%  '
$ws.Range('E43').Value = '  -3.47%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = "'17.70"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.43%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").Value = "'0.0265"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.42%  '
$ws.Range("D46").Value = "'9.38"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.53%  '
$ws.Range("D47").Value = "'2.72"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.32%  '
$ws.Range('D48').Value = '2.559.96'
$ws.Range('E48').Value = '  -4.24%  '
$ws.Range("D49").Value = "'93.11"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.01%  '
$ws.Range("D50").Value = "'70.74"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.19%  '
$ws.Range("D51").Value = "'50.71"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.25%  '
